$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")

# Add new config row: dbdriver / com.microsoft.sqlserver.jdbc.SQLServerDriver
$ws.Range("A8").Value = "dbdriver"
$ws.Range("B8").Value = "com.microsoft.sqlserver.jdbc.SQLServerDriver"

# Widen column B to fit the new longer value, keep C:D as before
# (ColumnWidth uses character units that Excel converts to the internal
#  width units; 41.29 round-trips to an internal column width of 42)
$ws.Columns.Item(2).ColumnWidth = 41.29

# Update the selection to B3 (matches new active cell in the diff)
$ws.Range("B3").Select()

$wb.Save()
